$wb = $excel.ActiveWorkbook

# Rename the existing "Add Book" sheet to "Book Details"
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Book Details"

# Add a new "Delete Book" sheet right after "Book Details"
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Delete Book"

# Populate the ID column header and values
$ws2.Range("A1").Value = "ID"
$ws2.Range("A2").Value = "abc1245"
$ws2.Range("A3").Value = "def2378"
$ws2.Range("A4").Value = "ghi4521"
$ws2.Range("A5").Value = "jkl3690"
$ws2.Range("A6").Value = "mno5987"

# Highlight the header cell with a yellow fill
$ws2.Range("A1").Interior.Color = 65535
